# Updated cryptos list with latest price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to remain text (avoids Excel auto-parsing
    # dotted numeric-looking strings like "327.79" into a number),
    # then restore the default "Normal" style so no formatting residue
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.908.86"
$ws.Range("E2").Value = "  -5.62%  "
Set-TextValue $ws.Range("D3") "1.822.08"
$ws.Range("E3").Value = "  -4.28%  "
$ws.Range("E4").Value = "  -0.36%  "
Set-TextValue $ws.Range("D5") "327.79"
$ws.Range("E5").Value = "  -3.13%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.27%  "
Set-TextValue $ws.Range("D7") "0.4633"
$ws.Range("E7").Value = "  -2.77%  "
Set-TextValue $ws.Range("D8") "0.3844"
$ws.Range("E8").Value = "  -3.92%  "
Set-TextValue $ws.Range("D9") "45.84"
$ws.Range("E9").Value = "  -2.82%  "
Set-TextValue $ws.Range("D10") "0.07850"
$ws.Range("E10").Value = "  -2.50%  "
Set-TextValue $ws.Range("D11") "0.9588"
$ws.Range("E11").Value = "  -3.30%  "
Set-TextValue $ws.Range("D12") "21.82"
$ws.Range("E12").Value = "  -5.89%  "
Set-TextValue $ws.Range("D13") "1.820.26"
$ws.Range("E13").Value = "  -3.37%  "
Set-TextValue $ws.Range("D14") "5.644"
$ws.Range("E14").Value = "  -4.82%  "
Set-TextValue $ws.Range("D16") "0.06868"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("E17").Value = "  -0.33%  "
Set-TextValue $ws.Range("D18") "86.38"
$ws.Range("E18").Value = "  -3.14%  "
Set-TextValue $ws.Range("D19") "0.000009923"
$ws.Range("E20").Value = "  -4.08%  "
Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  -0.38%  "
Set-TextValue $ws.Range("D22") "27.926.04"
$ws.Range("E22").Value = "  -5.55%  "
Set-TextValue $ws.Range("D23") "5.301"
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("E24").Value = "  -6.33%  "
$ws.Range("E25").Value = "  -2.59%  "
Set-TextValue $ws.Range("D26") "2.043.87"
$ws.Range("E26").Value = "  -3.22%  "
Set-TextValue $ws.Range("D27") "151.96"
$ws.Range("E27").Value = "  -3.00%  "
Set-TextValue $ws.Range("D28") "19.16"
$ws.Range("E28").Value = "  -2.33%  "
Set-TextValue $ws.Range("D29") "5.763"
$ws.Range("E29").Value = "  -11.48%  "
$ws.Range("E30").Value = "  -4.53%  "
Set-TextValue $ws.Range("D31") "116.59"
$ws.Range("E31").Value = "  -2.11%  "
Set-TextValue $ws.Range("D32") "0.9335"
$ws.Range("E32").Value = "  -6.32%  "
Set-TextValue $ws.Range("D33") "0.09215"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("E34").Value = "  -3.56%  "
Set-TextValue $ws.Range("D35") "1.313"
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("E36").Value = "  -5.33%  "
Set-TextValue $ws.Range("D37") "0.05928"
$ws.Range("E37").Value = "  -8.29%  "
$ws.Range("E38").Value = "  -4.46%  "
Set-TextValue $ws.Range("D39") "1.140"
$ws.Range("E39").Value = "  -4.67%  "
Set-TextValue $ws.Range("D40") "1.000"
$ws.Range("E40").Value = "  -0.49%  "
Set-TextValue $ws.Range("D41") "7.537"
$ws.Range("E41").Value = "  -2.43%  "
Set-TextValue $ws.Range("D42") "0.5567"
$ws.Range("E42").Value = "  -4.38%  "
Set-TextValue $ws.Range("D43") "9.926"
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("E44").Value = "  -3.20%  "
Set-TextValue $ws.Range("D45") "1.233"
$ws.Range("E45").Value = "  -2.59%  "
Set-TextValue $ws.Range("D46") "2.205"
Set-TextValue $ws.Range("D47") "11.54"
$ws.Range("E47").Value = "  -4.58%  "
Set-TextValue $ws.Range("D48") "0.5237"
$ws.Range("E48").Value = "  -4.42%  "
Set-TextValue $ws.Range("D49") "0.06993"
$ws.Range("E49").Value = "  -5.46%  "
Set-TextValue $ws.Range("D50") "1.819"
$ws.Range("E50").Value = "  -7.11%  "
Set-TextValue $ws.Range("D51") "111.96"
$ws.Range("E51").Value = "  -3.36%  "
